$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 87
$ws.Cells.Item(6, 6).Value = 123
$ws.Cells.Item(7, 6).Value = 847
$ws.Cells.Item(9, 6).Value = 1040
$ws.Cells.Item(10, 6).Value = 142
$ws.Cells.Item(11, 6).Value = 1069
$ws.Cells.Item(12, 6).Value = 806
$ws.Cells.Item(14, 6).Value = 677
$ws.Cells.Item(15, 6).Value = 1309
$ws.Cells.Item(16, 6).Value = 1026
$ws.Cells.Item(18, 6).Value = 746
$ws.Cells.Item(19, 6).Value = 730
$ws.Cells.Item(24, 6).Value = 1220
$ws.Cells.Item(25, 6).Value = 140
$ws.Cells.Item(26, 6).Value = 429
$ws.Cells.Item(28, 6).Value = 5136
$ws.Cells.Item(29, 6).Value = 245
$ws.Cells.Item(31, 6).Value = 2419
$ws.Cells.Item(32, 6).Value = 5801
$ws.Cells.Item(35, 6).Value = 588
$ws.Cells.Item(37, 6).Value = 1040
$ws.Cells.Item(39, 6).Value = 22
$ws.Cells.Item(41, 6).Value = 663

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(5, 6).Value = 2086
$ws.Cells.Item(15, 6).Value = 658
$ws.Cells.Item(41, 6).Value = 482

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 361

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(6, 6).Value = 87
$ws.Cells.Item(7, 6).Value = 361
$ws.Cells.Item(12, 6).Value = 847
$ws.Cells.Item(15, 6).Value = 1040
$ws.Cells.Item(16, 6).Value = 142
$ws.Cells.Item(17, 6).Value = 1069
$ws.Cells.Item(18, 6).Value = 806
$ws.Cells.Item(21, 6).Value = 677
$ws.Cells.Item(22, 6).Value = 1309
$ws.Cells.Item(24, 6).Value = 1026
$ws.Cells.Item(25, 6).Value = 746
$ws.Cells.Item(27, 6).Value = 730
$ws.Cells.Item(31, 6).Value = 1220
$ws.Cells.Item(32, 6).Value = 140
$ws.Cells.Item(33, 6).Value = 429
$ws.Cells.Item(35, 6).Value = 5136
$ws.Cells.Item(36, 6).Value = 245
$ws.Cells.Item(38, 6).Value = 2419
$ws.Cells.Item(39, 6).Value = 5801
$ws.Cells.Item(42, 6).Value = 588
$ws.Cells.Item(44, 6).Value = 1040
$ws.Cells.Item(45, 6).Value = 22
$ws.Cells.Item(46, 6).Value = 663
$ws.Cells.Item(50, 6).Value = 482
